# zk-emotion-proof.py.docx edit:
#  1. Strip the leftover "_dx_frag_StartFragment"/"_dx_frag_EndFragment"
#     clipboard bookmarks from the first paragraph.
#  2. Drop the unused `add` import.
#  3. Move the `print("ZK Proof Generated ... Verified")` line so it runs
#     right before `return proof` inside generate_proof(), instead of
#     after the final `proof = generate_proof(tile_800)` call.

$d = $word.ActiveDocument

# 1. Remove the hidden fragment bookmarks (structure-only, no visible text).
foreach ($name in @("_dx_frag_StartFragment", "_dx_frag_EndFragment")) {
    if ($d.Bookmarks.Exists($name)) {
        $d.Bookmarks.Item($name).Delete()
    }
}

# 2. from py_ecc.bn128 import G1, multiply, add  ->  ... import G1, multiply
$d.Content.Find.Execute("from py_ecc.bn128 import G1, multiply, add", $true, $false, $false, $false, $false,
                         $true, 1, $false, "from py_ecc.bn128 import G1, multiply", 2)

# 3. Relocate the print(...) line: insert it as a new paragraph right
#    before "    return proof", then delete the old trailing copy.
$printText = 'print("ZK Proof Generated — Fidelity 1.00 Verified")'

$returnPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd() -eq "    return proof") {
        $returnPara = $p
        break
    }
}

# InsertParagraphBefore() splices a new empty paragraph in immediately
# before the "return proof" text; the $returnPara handle then refers to
# that new (still empty) paragraph, so we can fill it in directly.
$returnPara.Range.InsertParagraphBefore()
$returnPara.Range.Text = "    " + $printText

# There are now two paragraphs with this text (the freshly inserted one
# and the original trailing one) - keep the first (new) occurrence and
# remove the last (original, trailing) occurrence.
$oldPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd() -eq $printText) {
        $oldPara = $p
    }
}
$oldPara.Range.Delete()
